# Update the "insurance" (保險, sheet index 5) and "debt" (債務, sheet index 6)
# sheets so that they follow the same column schema used by the other
# property-type sheets: a real header row (property_category, category,
# date, legislator_name, legislator_id, source_file, index, ...) instead of
# the header row duplicating row 2's data, plus the extra metadata columns
# filled in for every data row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "保險" (insurance) -- 5th tab
# ---------------------------------------------------------------------
$ins = $wb.Worksheets.Item(5)

# The "date" column (G) holds an ISO-looking string ("2012-04-25"); format
# it as text up front so it is stored as a literal string, not auto-parsed
# into a date serial number.
$ins.Range("G1:G4").NumberFormat = "@"

# Header row (row 1): B..K
$ins.Range("B1").Value = "company"
$ins.Range("C1").Value = "name"
$ins.Range("D1").Value = "owner"
$ins.Range("E1").Value = "property_category"
$ins.Range("F1").Value = "category"
$ins.Range("G1").Value = "date"
$ins.Range("H1").Value = "legislator_name"
$ins.Range("I1").Value = "legislator_id"
$ins.Range("J1").Value = "source_file"
$ins.Range("K1").Value = "index"

# Row 2 (index 125)
$ins.Range("B2").Value = "南山人壽"
$ins.Range("C2").Value = "南山金福利21年期養老壽險"
$ins.Range("D2").Value = "陳碧涵"
$ins.Range("E2").Value = "insurance"
$ins.Range("F2").Value = "normal"
$ins.Range("G2").Value = "2012-04-25"
$ins.Range("H2").Value = "陳碧涵"
$ins.Range("I2").Value = 1752
$ins.Range("J2").Value = "tmpd4df1"
$ins.Range("K2").Value = 125

# Row 3 (index 126)
$ins.Range("B3").Value = "台銀人壽股份有限公司"
$ins.Range("C3").Value = "鴻福還本終身壽險"
$ins.Range("D3").Value = "陳碧涵"
$ins.Range("E3").Value = "insurance"
$ins.Range("F3").Value = "normal"
$ins.Range("G3").Value = "2012-04-25"
$ins.Range("H3").Value = "陳碧涵"
$ins.Range("I3").Value = 1752
$ins.Range("J3").Value = "tmpd4df1"
$ins.Range("K3").Value = 126

# Row 4 (index 127)
$ins.Range("B4").Value = "南山人壽"
$ins.Range("C4").Value = "南山新康祥終身壽險"
$ins.Range("D4").Value = "廖◦陽"
$ins.Range("E4").Value = "insurance"
$ins.Range("F4").Value = "normal"
$ins.Range("G4").Value = "2012-04-25"
$ins.Range("H4").Value = "陳碧涵"
$ins.Range("I4").Value = 1752
$ins.Range("J4").Value = "tmpd4df1"
$ins.Range("K4").Value = 127

# ---------------------------------------------------------------------
# Sheet "債務" (debt) -- 6th tab
# ---------------------------------------------------------------------
$debt = $wb.Worksheets.Item(6)

# The "date" column (J) holds the ISO-looking string ("2012-04-25") too.
$debt.Range("J1:J6").NumberFormat = "@"

# Header row (row 1): B..N
$debt.Range("B1").Value = "species"
$debt.Range("C1").Value = "debtor"
$debt.Range("D1").Value = "owner"
$debt.Range("E1").Value = "total"
$debt.Range("F1").Value = "register_date"
$debt.Range("G1").Value = "register_reason"
$debt.Range("H1").Value = "property_category"
$debt.Range("I1").Value = "category"
$debt.Range("J1").Value = "date"
$debt.Range("K1").Value = "legislator_name"
$debt.Range("L1").Value = "legislator_id"
$debt.Range("M1").Value = "source_file"
$debt.Range("N1").Value = "index"

# Row 2 (index 137)
$debt.Range("B2").Value = "公教購屋貸款"
$debt.Range("C2").Value = "陳碧涵"
$debt.Range("D2").Value = "台灣銀行(水湳分行）臺中市北屯區祟德路"
$debt.Range("E2").Value = 2020088
$debt.Range("F2").Value = "94年12月13曰"
$debt.Range("G2").Value = "購屋貸款"
$debt.Range("H2").Value = "debt"
$debt.Range("I2").Value = "normal"
$debt.Range("J2").Value = "2012-04-25"
$debt.Range("K2").Value = "陳碧涵"
$debt.Range("L2").Value = 1752
$debt.Range("M2").Value = "tmpd4df1"
$debt.Range("N2").Value = 137

# Row 3 (index 138)
$debt.Range("B3").Value = "房屋抵押貸款"
$debt.Range("C3").Value = "陳碧涵"
$debt.Range("D3").Value = "台灣銀行(水湳分行）臺中市北屯區祟德路"
$debt.Range("E3").Value = 2747491
$debt.Range("F3").Value = "94年12月13日"
$debt.Range("G3").Value = "個人資金調度"
$debt.Range("H3").Value = "debt"
$debt.Range("I3").Value = "normal"
$debt.Range("J3").Value = "2012-04-25"
$debt.Range("K3").Value = "陳碧涵"
$debt.Range("L3").Value = 1752
$debt.Range("M3").Value = "tmpd4df1"
$debt.Range("N3").Value = 138

# Row 4 (index 139)
$debt.Range("B4").Value = "房屋抵押貸款"
$debt.Range("C4").Value = "陳碧涵"
$debt.Range("D4").Value = "第一銀行(進化分行）臺中市北屯區進化北路"
$debt.Range("E4").Value = 460745
$debt.Range("F4").Value = "96年08月30日"
$debt.Range("G4").Value = "個人資金調度"
$debt.Range("H4").Value = "debt"
$debt.Range("I4").Value = "normal"
$debt.Range("J4").Value = "2012-04-25"
$debt.Range("K4").Value = "陳碧涵"
$debt.Range("L4").Value = 1752
$debt.Range("M4").Value = "tmpd4df1"
$debt.Range("N4").Value = 139

# Row 5 (index 141)
$debt.Range("B5").Value = "房屋抵押貸款"
$debt.Range("C5").Value = "陳碧涵"
$debt.Range("D5").Value = "第一銀行(進化分行）臺中市北屯區進化北路"
$debt.Range("E5").Value = 1213710
$debt.Range("F5").Value = "89年01月27日"
$debt.Range("G5").Value = "個人資金調度"
$debt.Range("H5").Value = "debt"
$debt.Range("I5").Value = "normal"
$debt.Range("J5").Value = "2012-04-25"
$debt.Range("K5").Value = "陳碧涵"
$debt.Range("L5").Value = 1752
$debt.Range("M5").Value = "tmpd4df1"
$debt.Range("N5").Value = 141

# Row 6 (index 142)
$debt.Range("B6").Value = "信用貸款"
$debt.Range("C6").Value = "陳碧涵"
$debt.Range("D6").Value = "合作金庫(中興分行）臺中市北屯區公圜路"
$debt.Range("E6").Value = 15436
$debt.Range("F6").Value = "94年10月18闩"
$debt.Range("G6").Value = "個人資金調度"
$debt.Range("H6").Value = "debt"
$debt.Range("I6").Value = "normal"
$debt.Range("J6").Value = "2012-04-25"
$debt.Range("K6").Value = "陳碧涵"
$debt.Range("L6").Value = 1752
$debt.Range("M6").Value = "tmpd4df1"
$debt.Range("N6").Value = 142
